# Update "想去人数" (want-to-go count) figures in column F on the
# "展览" and "全部类型" sheets, reflecting the refreshed scrape output.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 5308
$ws1.Range("F22").Value = 3489
$ws1.Range("F23").Value = 1094
$ws1.Range("F24").Value = 2776
$ws1.Range("F27").Value = 3996
$ws1.Range("F32").Value = 24
$ws1.Range("F35").Value = 54
$ws1.Range("F36").Value = 1012
$ws1.Range("F38").Value = 499
$ws1.Range("F41").Value = 3550

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 5308
$ws4.Range("F24").Value = 3489
$ws4.Range("F27").Value = 1094
$ws4.Range("F29").Value = 2776
$ws4.Range("F31").Value = 3996
$ws4.Range("F36").Value = 24
$ws4.Range("F40").Value = 54
$ws4.Range("F41").Value = 1012
$ws4.Range("F48").Value = 3550
